$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# Sheet "OR stunting by compfeeding": add 4 new rows (7-10) that mirror
# rows 2-5 but for the new "IYCF" complementary-feeding categories.
# ------------------------------------------------------------------
$ws17 = $wb.Worksheets.Item("OR stunting by compfeeding")

# Clone formatting of the existing data block down into the new rows,
# then set row heights so the new rows get the same "custom height" row
# metadata as the rest of the sheet.
$ws17.Range("A2:F5").Copy()
$ws17.Range("A7").PasteSpecial(-4122)
$ws17.Rows.Item(7).RowHeight = 15.75
$ws17.Rows.Item(8).RowHeight = 15.75
$ws17.Rows.Item(9).RowHeight = 15.75
$ws17.Rows.Item(10).RowHeight = 15.75

$ws17.Range("A7").Value = "Complementary feeding (food secure with IYCF)"
$ws17.Range("B7").Value = 1
$ws17.Range("C7").Value = 1
$ws17.Range("D7").Value = 1
$ws17.Range("E7").Value = 1
$ws17.Range("F7").Value = 1

$ws17.Range("A8").Value = "Complementary feeding (food secure without IYCF)"
$ws17.Range("B8").Value = 1
$ws17.Range("C8").Value = 1
$ws17.Range("D8").Value = 1.43
$ws17.Range("E8").Value = 1.43
$ws17.Range("F8").Value = 1

$ws17.Range("A9").Value = "Complementary feeding (food insecure with IYCF and supplementation)"
$ws17.Range("B9").Value = 1
$ws17.Range("C9").Value = 1
$ws17.Range("D9").Value = 1.6
$ws17.Range("E9").Value = 1.6
$ws17.Range("F9").Value = 1

$ws17.Range("A10").Value = "Complementary feeding (food insecure with neither IYCF nor supplementation)"
$ws17.Range("B10").Value = 1
$ws17.Range("C10").Value = 1
$ws17.Range("D10").Value = 2.39
$ws17.Range("E10").Value = 2.39
$ws17.Range("F10").Value = 1

$ws17.Range("A7:F10").Select()

# ------------------------------------------------------------------
# Sheet "OR correctBF by interventn": fill in the previously-blank row 4
# with a new "IYCF" intervention row (mirrors row 3's formatting).
# ------------------------------------------------------------------
$ws18 = $wb.Worksheets.Item("OR correctBF by interventn")

$ws18.Range("A3:F3").Copy()
$ws18.Range("A4").PasteSpecial(-4122)

$ws18.Range("A4").Value = "IYCF"
$ws18.Range("B4").Value = 5.16
$ws18.Range("C4").Value = 5.16
$ws18.Range("D4").Value = 1.82
$ws18.Range("E4").Value = 1.82
$ws18.Range("F4").Value = 1

$ws18.Range("A4:F4").Select()

# ------------------------------------------------------------------
# Sheet "Interventions cost and coverage": add a new "IYCF" row (row 9),
# highlighted with a new peach fill colour, mirroring row 8's layout.
# Apply the fill to the fontId-0-styled cell (C9) first and the
# fontId-4-styled cells (B9/D9) second so the new style entries land in
# the same order as the source workbook.
# ------------------------------------------------------------------
$ws20 = $wb.Worksheets.Item("Interventions cost and coverage")

$ws20.Range("C8").Copy()
$ws20.Range("C9").PasteSpecial(-4122)
$ws20.Range("C9").Interior.Color = 14281213

$ws20.Range("B8").Copy()
$ws20.Range("B9").PasteSpecial(-4122)
$ws20.Range("B9").Interior.Color = 14281213

$ws20.Range("D8").Copy()
$ws20.Range("D9").PasteSpecial(-4122)
$ws20.Range("D9").Interior.Color = 14281213

$ws20.Range("A8").Copy()
$ws20.Range("A9").PasteSpecial(-4122)

$ws20.Range("A9").Value = "IYCF"
$ws20.Range("B9").Value = 0
$ws20.Range("C9").Value = 0.95
$ws20.Range("D9").Value = 10.49

# Activating this sheet last makes it the active tab (matches the
# workbook's new activeTab / tabSelected state) and moves the selection
# here, which also clears tabSelected from whichever sheet had it before.
$ws20.Activate()
$ws20.Range("A9:D9").Select()
